$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = 406600
$ws.Range("E8").Value = 313000
$ws.Range("F8").Value = 215300
$ws.Range("G8").Value = 126300
$ws.Range("H8").Value = 84100
$ws.Range("I8").Value = 66800
$ws.Range("D9").Value = 279100
$ws.Range("E9").Value = 224900
$ws.Range("F9").Value = 168900
$ws.Range("G9").Value = 106700
$ws.Range("H9").Value = 78100
$ws.Range("I9").Value = 64200
$ws.Range("D10").Value = 127500
$ws.Range("E10").Value = 88100
$ws.Range("F10").Value = 46400
$ws.Range("G10").Value = 19700
$ws.Range("H10").Value = 5900
$ws.Range("D14").Value = 3000
$ws.Range("D17").Value = 342300
$ws.Range("E17").Value = 275200
$ws.Range("F17").Value = 204200
$ws.Range("G17").Value = 129000
$ws.Range("H17").Value = 98800
$ws.Range("I17").Value = 82200
$ws.Range("D18").Value = 64300
$ws.Range("E18").Value = 37800
$ws.Range("F18").Value = 11100
$ws.Range("H18").Value = -14700
$ws.Range("I18").Value = -15400
$ws.Range("E20").Value = 1500
$ws.Range("F20").Value = 123600
$ws.Range("D21").Value = 179100
$ws.Range("E21").Value = 138600
$ws.Range("F21").Value = 205300
$ws.Range("G21").Value = 40600
$ws.Range("H21").Value = 14300
$ws.Range("I21").Value = 9400
$ws.Range("D22").Value = 41600
$ws.Range("E22").Value = 33400
$ws.Range("F22").Value = 18400
$ws.Range("G22").Value = 11400
$ws.Range("H22").Value = 7600
$ws.Range("I22").Value = 9900
$ws.Range("D23").Value = 25200
$ws.Range("E23").Value = 5900
$ws.Range("F23").Value = 116300
$ws.Range("G23").Value = -13500
$ws.Range("H23").Value = -22400
$ws.Range("I23").Value = -25300
$ws.Range("D24").Value = 7000
$ws.Range("F24").Value = 13000
$ws.Range("D26").Value = 18100
$ws.Range("E26").Value = 4900
$ws.Range("F26").Value = 103300
$ws.Range("G26").Value = -13800
$ws.Range("H26").Value = -22600
$ws.Range("I26").Value = -26100
$ws.Range("D27").Value = 18100
$ws.Range("E27").Value = 4900
$ws.Range("F27").Value = 103300
$ws.Range("G27").Value = -51000
$ws.Range("H27").Value = -55200
$ws.Range("I27").Value = -49100
$ws.Range("E32").Value = -1500
$ws.Range("F32").Value = -123600
$ws.Range("D33").Value = 18100
$ws.Range("E33").Value = 4900
$ws.Range("F33").Value = 103300
$ws.Range("G33").Value = -51000
$ws.Range("H33").Value = -55200
$ws.Range("I33").Value = -49100
$ws.Range("D35").Value = 18100
$ws.Range("E35").Value = 4900
$ws.Range("F35").Value = 103300
$ws.Range("G35").Value = -51000
$ws.Range("H35").Value = -55200
$ws.Range("I35").Value = -49100
$ws.Range("D41").Value = 99600
$ws.Range("E41").Value = 78600
$ws.Range("F41").Value = 387400
$ws.Range("G41").Value = 137500
$ws.Range("H41").Value = 93600
$ws.Range("I41").Value = 19800
$ws.Range("D43").Value = 240600
$ws.Range("E43").Value = 136400
$ws.Range("F43").Value = 55400
$ws.Range("G43").Value = 28300
$ws.Range("H43").Value = 10900
$ws.Range("I43").Value = 8200
$ws.Range("D45").Value = 159100
$ws.Range("E45").Value = 213100
$ws.Range("F45").Value = 65800
$ws.Range("G45").Value = 45900
$ws.Range("H45").Value = 14800
$ws.Range("I45").Value = 7500
$ws.Range("D46").Value = 499300
$ws.Range("E46").Value = 288200
$ws.Range("F46").Value = 508600
$ws.Range("G46").Value = 211700
$ws.Range("H46").Value = 119300
$ws.Range("I46").Value = 35500
$ws.Range("G47").Value = 22700
$ws.Range("D48").Value = 1047300
$ws.Range("E48").Value = 1698900
$ws.Range("F48").Value = 608000
$ws.Range("G48").Value = 287900
$ws.Range("H48").Value = 157700
$ws.Range("I48").Value = 125300
$ws.Range("D49").Value = 9900
$ws.Range("E49").Value = 9500
$ws.Range("F49").Value = 6700
$ws.Range("G49").Value = 5700
$ws.Range("I49").Value = 4300
$ws.Range("D52").Value = 63800
$ws.Range("E52").Value = 64100
$ws.Range("F52").Value = 34400
$ws.Range("G52").Value = 29400
$ws.Range("H52").Value = 19300
$ws.Range("D54").Value = 1621900
$ws.Range("E54").Value = 1211200
$ws.Range("F54").Value = 1157700
$ws.Range("G54").Value = 557400
$ws.Range("H54").Value = 300700
$ws.Range("I54").Value = 165700
$ws.Range("D57").Value = 114600
$ws.Range("E57").Value = 77500
$ws.Range("F57").Value = 116600
$ws.Range("D58").Value = 394400
$ws.Range("E58").Value = 274900
$ws.Range("F58").Value = 119200
$ws.Range("G58").Value = 80200
$ws.Range("H58").Value = 32600
$ws.Range("I58").Value = 60800
$ws.Range("D59").Value = 44200
$ws.Range("E59").Value = 34400
$ws.Range("F59").Value = 43500
$ws.Range("G59").Value = 19300
$ws.Range("H59").Value = 15900
$ws.Range("I59").Value = 17300
$ws.Range("D60").Value = 553200
$ws.Range("E60").Value = 207200
$ws.Range("F60").Value = 279300
$ws.Range("G60").Value = 100300
$ws.Range("H60").Value = 49500
$ws.Range("I60").Value = 78900
$ws.Range("D61").Value = 438900
$ws.Range("E61").Value = 410800
$ws.Range("F61").Value = 292300
$ws.Range("G61").Value = 105900
$ws.Range("H61").Value = 55800
$ws.Range("I61").Value = 1000
$ws.Range("D62").Value = 5500
$ws.Range("H70").Value = 337400
$ws.Range("I70").Value = 173600
$ws.Range("D72").Value = -60200
$ws.Range("E72").Value = -78300
$ws.Range("F72").Value = -83200
$ws.Range("G72").Value = -186600
$ws.Range("H72").Value = -143000
$ws.Range("I72").Value = -88900
$ws.Range("D76").Value = 624300
$ws.Range("E76").Value = 592300
$ws.Range("F76").Value = 585900
$ws.Range("G76").Value = 351200
$ws.Range("H76").Value = -142000
$ws.Range("I76").Value = -88500
$ws.Range("D81").Value = 18100
$ws.Range("E81").Value = 4900
$ws.Range("F81").Value = 103300
$ws.Range("G81").Value = -51000
$ws.Range("H81").Value = -55200
$ws.Range("I81").Value = -49100
$ws.Range("D83").Value = 112100
$ws.Range("E83").Value = 99100
$ws.Range("F83").Value = 70500
$ws.Range("G83").Value = 42700
$ws.Range("H83").Value = 29100
$ws.Range("I83").Value = 24800
$ws.Range("D89").Value = 106800
$ws.Range("E89").Value = 62400
$ws.Range("F89").Value = 44000
$ws.Range("G89").Value = 6400
$ws.Range("I89").Value = 3900
$ws.Range("D91").Value = -482700
$ws.Range("E91").Value = -547300
$ws.Range("F91").Value = -324000
$ws.Range("G91").Value = -197000
$ws.Range("H91").Value = -89200
$ws.Range("D94").Value = -383700
$ws.Range("E94").Value = -501900
$ws.Range("F94").Value = -147200
$ws.Range("G94").Value = -231800
$ws.Range("H94").Value = -85600
$ws.Range("I94").Value = -40500
$ws.Range("D100").Value = 301600
$ws.Range("E100").Value = 118500
$ws.Range("F100").Value = 338600
$ws.Range("G100").Value = 269200
$ws.Range("H100").Value = 157800
$ws.Range("I100").Value = 47700
$ws.Range("E101").Value = 12200
$ws.Range("F101").Value = 14500
$ws.Range("D102").Value = 21000
$ws.Range("E102").Value = -308800
$ws.Range("F102").Value = 249900
$ws.Range("G102").Value = 43900
$ws.Range("H102").Value = 73800
$ws.Range("I102").Value = 11100
